$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values (coin prices) are stored as text, matching the
# original inline-string cells, so numeric-looking strings like "210.62"
# are not silently reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.949.63'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = '1.595.65'
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("E4").Value = '  +0.74%  '
$ws.Range("D5").Value = '210.62'
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("D7").Value = '0.482'
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("D9").Value = '0.0610'
$ws.Range("E9").Value = '  -1.70%  '
$ws.Range("D10").Value = '18.04'
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("D11").Value = '0.0809'
$ws.Range("E11").Value = '  +3.30%  '
$ws.Range("D12").Value = '1.820.14'
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("D13").Value = '1.600.93'
$ws.Range("D14").Value = '3.99'
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D16").Value = '25.967.89'
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("D17").Value = '59.98'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("D20").Value = '200.11'
$ws.Range("E20").Value = '  +4.76%  '
$ws.Range("D21").Value = '4.22'
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").Value = '9.24'
$ws.Range("E22").Value = '  -2.27%  '
$ws.Range("E23").Value = '  +0.68%  '
$ws.Range("E24").Value = '  +5.00%  '
$ws.Range("D25").Value = '141.68'
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("E27").Value = '  -8.71%  '
$ws.Range("D28").Value = '15.08'
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("E33").Value = '  -2.75%  '
$ws.Range("D34").Value = '1.47'
$ws.Range("E34").Value = '  -1.87%  '
$ws.Range("E35").Value = '  +2.79%  '
$ws.Range("D36").Value = '1.123.37'
$ws.Range("E36").Value = '  +1.84%  '
$ws.Range("E37").Value = '  +7.36%  '
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("E39").Value = '  -0.74%  '
$ws.Range("D40").Value = '0.784'
$ws.Range("E40").Value = '  -0.98%  '
$ws.Range("D41").Value = '0.489'
$ws.Range("E41").Value = '  -3.21%  '
$ws.Range("D42").Value = '0.779'
$ws.Range("E42").Value = '  -3.61%  '
$ws.Range("D43").Value = '1.731.28'
$ws.Range("E43").Value = '  +1.56%  '
$ws.Range("D44").Value = '5.10'
$ws.Range("E44").Value = '  -1.00%  '
$ws.Range("D45").Value = '92.60'
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("D47").Value = '53.27'
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("E48").Value = '  -1.20%  '
$ws.Range("D49").Value = '0.407'
$ws.Range("E49").Value = '  +0.69%  '
$ws.Range("D50").Value = '1.01'
$ws.Range("E50").Value = '  +0.86%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₇0918'
$ws.Range("E51").Value = '  -12.52%  '
